$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 11th week label into the next empty row of column A
$ws.Range("A11").Value = "FC-w11"
